$d = $word.ActiveDocument

# Locate the first empty paragraph immediately following the last
# "        }" code line (the 4 trailing empty paragraphs before sectPr).
$count = $d.Paragraphs.Count
$target = $null
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim().Length -gt 0) {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}

# Center-align the (previously empty) paragraph and move the "_GoBack"
# bookmark here (Word keeps a single "_GoBack" bookmark, so adding it
# here removes it from its old location automatically).
$target.Format.Alignment = 1
$d.Bookmarks.Add("_GoBack", $target.Range)
